# Update countries & provincias Spain
# - Refresh COVID-19 case counters for several countries (rows shown below)
# - Insert "Malaui" ahead of "Togo" in the country list (rows 144-149 shift down one slot)
# - Reorder a small cluster of territories/islands (rows 206-214)
# - Bump the "Datos actualizados" timestamp to 19:51

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - Datos actualizados a 13 de Junio de 2020 a las 19:51
$ws.Range("A1").Value = 'Datos actualizados a 13 de Junio de 2020 a las 19:51'

# Row 4 - Estados Unidos
$ws.Range("A4").Value = 'Estados Unidos'
$ws.Range("B4").Value = 2129736
$ws.Range("C4").Value = 12814
$ws.Range("D4").Value = 843548
$ws.Range("E4").Value = 1169092
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 271
$ws.Range("H4").Value = 117096

# Row 7 - India
$ws.Range("A7").Value = 'India'
$ws.Range("B7").Value = 321406
$ws.Range("C7").Value = 11803
$ws.Range("D7").Value = 162320
$ws.Range("E7").Value = 149881
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 315
$ws.Range("H7").Value = 9205

# Row 9 - España
$ws.Range("A9").Value = 'España'
$ws.Range("B9").Value = 290685
$ws.Range("C9").Value = 396
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 27136

# Row 12 - Alemania
$ws.Range("A12").Value = 'Alemania'
$ws.Range("B12").Value = 187356
$ws.Range("C12").Value = 105
$ws.Range("D12").Value = 171900
$ws.Range("E12").Value = 6592
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = 8864

# Row 14 - Turquia
$ws.Range("A14").Value = 'Turquia'
$ws.Range("B14").Value = 176677
$ws.Range("C14").Value = 1459
$ws.Range("D14").Value = 150087
$ws.Range("E14").Value = 21798
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 14
$ws.Range("H14").Value = 4792

# Row 21 - Banglades
$ws.Range("A21").Value = 'Banglades'
$ws.Range("B21").Value = 84379
$ws.Range("C21").Value = 2856
$ws.Range("D21").Value = 17828
$ws.Range("E21").Value = 65412
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 44
$ws.Range("H21").Value = 1139

# Row 42 - Irlanda
$ws.Range("A42").Value = 'Irlanda'
$ws.Range("B42").Value = 25295
$ws.Range("C42").Value = 45
$ws.Range("D42").Value = 22698
$ws.Range("E42").Value = 892
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = 1705

# Row 66 - Marruecos
$ws.Range("A66").Value = 'Marruecos'
$ws.Range("B66").Value = 8692
$ws.Range("C66").Value = 82
$ws.Range("D66").Value = 7696
$ws.Range("E66").Value = 784
$ws.Range("F66").Value = 0
$ws.Range("G66").Value = 0
$ws.Range("H66").Value = 212

# Row 101 - Maldivas
$ws.Range("A101").Value = 'Maldivas'
$ws.Range("B101").Value = 2013
$ws.Range("C101").Value = 10
$ws.Range("D101").Value = 1217
$ws.Range("E101").Value = 788
$ws.Range("F101").Value = 0
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 8

# Row 119 - Paraguay
$ws.Range("A119").Value = 'Paraguay'
$ws.Range("B119").Value = 1261
$ws.Range("C119").Value = 7
$ws.Range("D119").Value = 647
$ws.Range("E119").Value = 603
$ws.Range("F119").Value = 0
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 11

# Row 132 - Republica del Chad
$ws.Range("A132").Value = 'Republica del Chad'
$ws.Range("B132").Value = 848
$ws.Range("C132").Value = 0
$ws.Range("D132").Value = 718
$ws.Range("E132").Value = 58
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 0
$ws.Range("H132").Value = 72

# Row 144 - Malaui
$ws.Range("A144").Value = 'Malaui'
$ws.Range("B144").Value = 529
$ws.Range("C144").Value = 48
$ws.Range("D144").Value = 66
$ws.Range("E144").Value = 458
$ws.Range("F144").Value = 0
$ws.Range("G144").Value = 1
$ws.Range("H144").Value = 5

# Row 145 - Togo
$ws.Range("A145").Value = 'Togo'
$ws.Range("B145").Value = 525
$ws.Range("C145").Value = 0
$ws.Range("D145").Value = 279
$ws.Range("E145").Value = 233
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 13

# Row 146 - Ruanda
$ws.Range("A146").Value = 'Ruanda'
$ws.Range("B146").Value = 510
$ws.Range("C146").Value = 0
$ws.Range("D146").Value = 321
$ws.Range("E146").Value = 187
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 2

# Row 147 - Tanzania
$ws.Range("A147").Value = 'Tanzania'
$ws.Range("B147").Value = 509
$ws.Range("C147").Value = 0
$ws.Range("D147").Value = 183
$ws.Range("E147").Value = 305
$ws.Range("F147").Value = 0
$ws.Range("G147").Value = 0
$ws.Range("H147").Value = 21

# Row 148 - Estado de Palestina
$ws.Range("A148").Value = 'Estado de Palestina'
$ws.Range("B148").Value = 489
$ws.Range("C148").Value = 0
$ws.Range("D148").Value = 414
$ws.Range("E148").Value = 72
$ws.Range("F148").Value = 0
$ws.Range("G148").Value = 0
$ws.Range("H148").Value = 3

# Row 149 - Reunion
$ws.Range("A149").Value = 'Reunion'
$ws.Range("B149").Value = 489
$ws.Range("C149").Value = 1
$ws.Range("D149").Value = 460
$ws.Range("E149").Value = 28
$ws.Range("F149").Value = 0
$ws.Range("G149").Value = 0
$ws.Range("H149").Value = 1

# Row 155 - Zimbabue
$ws.Range("A155").Value = 'Zimbabue'
$ws.Range("B155").Value = 356
$ws.Range("C155").Value = 13
$ws.Range("D155").Value = 54
$ws.Range("E155").Value = 298
$ws.Range("F155").Value = 0
$ws.Range("G155").Value = 0
$ws.Range("H155").Value = 4

# Row 206 - Islas Malvinas
$ws.Range("A206").Value = 'Islas Malvinas'
$ws.Range("B206").Value = 13
$ws.Range("C206").Value = 0
$ws.Range("D206").Value = 13
$ws.Range("E206").Value = 0
$ws.Range("F206").Value = 0
$ws.Range("G206").Value = 0
$ws.Range("H206").Value = 0

# Row 207 - Groenlandia
$ws.Range("A207").Value = 'Groenlandia'
$ws.Range("B207").Value = 13
$ws.Range("C207").Value = 0
$ws.Range("D207").Value = 13
$ws.Range("E207").Value = 0
$ws.Range("F207").Value = 0
$ws.Range("G207").Value = 0
$ws.Range("H207").Value = 0

# Row 208 - Islas Turcas y Caicos
$ws.Range("A208").Value = 'Islas Turcas y Caicos'
$ws.Range("B208").Value = 12
$ws.Range("C208").Value = 0
$ws.Range("D208").Value = 11
$ws.Range("E208").Value = 0
$ws.Range("F208").Value = 0
$ws.Range("G208").Value = 0
$ws.Range("H208").Value = 1

# Row 209 - Santa Sede
$ws.Range("A209").Value = 'Santa Sede'
$ws.Range("B209").Value = 12
$ws.Range("C209").Value = 0
$ws.Range("D209").Value = 12
$ws.Range("E209").Value = 0
$ws.Range("F209").Value = 0
$ws.Range("G209").Value = 0
$ws.Range("H209").Value = 0

# Row 210 - Seychelles
$ws.Range("A210").Value = 'Seychelles'
$ws.Range("B210").Value = 11
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 11
$ws.Range("E210").Value = 0
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 0

# Row 211 - Montserrat
$ws.Range("A211").Value = 'Montserrat'
$ws.Range("B211").Value = 11
$ws.Range("C211").Value = 0
$ws.Range("D211").Value = 10
$ws.Range("E211").Value = 0
$ws.Range("F211").Value = 0
$ws.Range("G211").Value = 0
$ws.Range("H211").Value = 1

# Row 213 - Papua Nueva Guinea
$ws.Range("A213").Value = 'Papua Nueva Guinea'
$ws.Range("B213").Value = 8
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 8
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 0

# Row 214 - Islas Virgenes Britanicas
$ws.Range("A214").Value = 'Islas Virgenes Britanicas'
$ws.Range("B214").Value = 8
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 7
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 1
